# Auto-generated edit script: updates profit-calculation cells in the
# Exodus_Profits workbook's per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market-board pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("H70").Value = 2074.5
$ws.Range("I70").Value = 2333
$ws.Range("J70").Value = 1919.4
$ws.Range("K70").Value = 6999
$ws.Range("L70").Value = 5758.200000000001
$ws.Range("M70").Value = -6729
$ws.Range("N70").Value = -6298.200000000001
$ws.Range("H73").Value = 2074.5
$ws.Range("I73").Value = 2333
$ws.Range("J73").Value = 1919.4
$ws.Range("K73").Value = 6999
$ws.Range("L73").Value = 5758.200000000001
$ws.Range("M73").Value = -6063
$ws.Range("N73").Value = -7630.200000000001
$ws.Range("H113").Value = 2995.4
$ws.Range("I113").Value = 2995.4
$ws.Range("K113").Value = 2995.4
$ws.Range("M113").Value = 258.5999999999999
$ws.Range("H114").Value = 44748.75
$ws.Range("J114").Value = 52998.668
$ws.Range("L114").Value = 52998.668
$ws.Range("N114").Value = -61676.668
$ws.Range("H117").Value = 56665
$ws.Range("J117").Value = 56665
$ws.Range("L117").Value = 56665
$ws.Range("N117").Value = -65843
$ws.Range("H134").Value = 84320.55499999999
$ws.Range("J134").Value = 84320.55499999999
$ws.Range("L134").Value = 84320.55499999999
$ws.Range("N134").Value = -94460.55499999999
$ws.Range("H136").Value = 77537.78
$ws.Range("J136").Value = 77537.78
$ws.Range("L136").Value = 77537.78
$ws.Range("N136").Value = -87737.78
$ws.Range("H139").Value = 69346.766
$ws.Range("J139").Value = 69346.766
$ws.Range("L139").Value = 69346.766
$ws.Range("N139").Value = -79626.766
$ws.Range("N44").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5644.356
$ws.Range("I32").Value = 2588.5293
$ws.Range("K32").Value = 2588.5293
$ws.Range("M32").Value = -2301.5293
$ws.Range("H52").Value = 46998
$ws.Range("J52").Value = 46998
$ws.Range("L52").Value = 46998
$ws.Range("N52").Value = -47634
$ws.Range("H110").Value = 1270
$ws.Range("I110").Value = 951
$ws.Range("J110").Value = 2333.3333
$ws.Range("K110").Value = 951
$ws.Range("L110").Value = 2333.3333
$ws.Range("M110").Value = 1094
$ws.Range("N110").Value = -6423.3333
$ws.Range("H121").Value = 71006.5
$ws.Range("J121").Value = 71006.5
$ws.Range("L121").Value = 71006.5
$ws.Range("N121").Value = -74500.5
$ws.Range("H122").Value = 2530.7144
$ws.Range("I122").Value = 2827.5
$ws.Range("K122").Value = 8482.5
$ws.Range("M122").Value = -6032.5
$ws.Range("H127").Value = 89082.86
$ws.Range("J127").Value = 89082.86
$ws.Range("L127").Value = 89082.86
$ws.Range("N127").Value = -99002.86

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 59882
$ws.Range("J13").Value = 59882
$ws.Range("L13").Value = 59882
$ws.Range("N13").Value = -60218
$ws.Range("H55").Value = 27731.166
$ws.Range("J55").Value = 27731.166
$ws.Range("L55").Value = 27731.166
$ws.Range("N55").Value = -28277.166
$ws.Range("H109").Value = 71997.14
$ws.Range("J109").Value = 71997.14
$ws.Range("L109").Value = 71997.14
$ws.Range("N109").Value = -74771.14
$ws.Range("H114").Value = 61494.332
$ws.Range("J114").Value = 61494.332
$ws.Range("L114").Value = 61494.332
$ws.Range("N114").Value = -70172.33199999999
$ws.Range("H115").Value = 76998
$ws.Range("J115").Value = 76998
$ws.Range("L115").Value = 76998
$ws.Range("N115").Value = -80132
$ws.Range("H117").Value = 87884.8
$ws.Range("J117").Value = 87884.8
$ws.Range("L117").Value = 87884.8
$ws.Range("N117").Value = -97062.8
$ws.Range("H135").Value = 97665.71000000001
$ws.Range("J135").Value = 97665.71000000001
$ws.Range("L135").Value = 97665.71000000001
$ws.Range("N135").Value = -107805.71
$ws.Range("H138").Value = 90862.22
$ws.Range("J138").Value = 90862.22
$ws.Range("L138").Value = 90862.22
$ws.Range("N138").Value = -101142.22
$ws.Range("H140").Value = 79995.664
$ws.Range("J140").Value = 79995.664
$ws.Range("L140").Value = 79995.664
$ws.Range("N140").Value = -90355.664

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 38490.5
$ws.Range("J18").Value = 38490.5
$ws.Range("L18").Value = 38490.5
$ws.Range("N18").Value = -38950.5
$ws.Range("H51").Value = 25733.334
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H59").Value = 121666
$ws.Range("J59").Value = 121666
$ws.Range("L59").Value = 121666
$ws.Range("N59").Value = -123956
$ws.Range("H61").Value = 25733.334
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("H62").Value = 2665.6667
$ws.Range("I62").Value = 2239.8
$ws.Range("K62").Value = 2239.8
$ws.Range("M62").Value = -1615.8
$ws.Range("H65").Value = 2665.6667
$ws.Range("I65").Value = 2239.8
$ws.Range("K65").Value = 11199
$ws.Range("M65").Value = -8079
$ws.Range("H86").Value = 2765861.8
$ws.Range("I86").Value = 5965204
$ws.Range("J86").Value = 23568.428
$ws.Range("K86").Value = 5965204
$ws.Range("L86").Value = 23568.428
$ws.Range("M86").Value = -5964081
$ws.Range("N86").Value = -25814.428
$ws.Range("H89").Value = 2765861.8
$ws.Range("I89").Value = 5965204
$ws.Range("J89").Value = 23568.428
$ws.Range("K89").Value = 29826020
$ws.Range("L89").Value = 117842.14
$ws.Range("M89").Value = -29820404
$ws.Range("N89").Value = -129074.14
$ws.Range("H118").Value = 173835.12
$ws.Range("J118").Value = 173835.12
$ws.Range("L118").Value = 173835.12
$ws.Range("N118").Value = -177149.12
$ws.Range("H122").Value = 5076.5557
$ws.Range("I122").Value = 4148.1665
$ws.Range("J122").Value = 6933.3335
$ws.Range("K122").Value = 12444.4995
$ws.Range("L122").Value = 20800.0005
$ws.Range("M122").Value = -9994.499500000002
$ws.Range("N122").Value = -25700.0005
$ws.Range("H138").Value = 99981.664
$ws.Range("J138").Value = 99981.664
$ws.Range("L138").Value = 99981.664
$ws.Range("N138").Value = -110261.664
$ws.Range("N51").ClearContents()
$ws.Range("N61").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4444663
$ws.Range("I4").Value = 4545650.5
$ws.Range("J4").Value = 4000320
$ws.Range("K4").Value = 13636951.5
$ws.Range("L4").Value = 12000960
$ws.Range("M4").Value = -13636839.5
$ws.Range("N4").Value = -12001184
$ws.Range("H48").Value = 3619.8
$ws.Range("I48").Value = 2449.5
$ws.Range("J48").Value = 4400
$ws.Range("K48").Value = 7348.5
$ws.Range("L48").Value = 13200
$ws.Range("M48").Value = -7098.5
$ws.Range("N48").Value = -13700

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 35665.668
$ws.Range("J109").Value = 35665.668
$ws.Range("L109").Value = 35665.668
$ws.Range("N109").Value = -37745.668
$ws.Range("H135").Value = 97183.81
$ws.Range("J135").Value = 97183.81
$ws.Range("L135").Value = 97183.81
$ws.Range("N135").Value = -107323.81
$ws.Range("H138").Value = 113333.336
$ws.Range("J138").Value = 113333.336
$ws.Range("L138").Value = 113333.336
$ws.Range("N138").Value = -123613.336
$ws.Range("H140").Value = 90251.75
$ws.Range("J140").Value = 90251.75
$ws.Range("L140").Value = 90251.75
$ws.Range("N140").Value = -100611.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6025.875
$ws.Range("J61").Value = 4949.5
$ws.Range("L61").Value = 4949.5
$ws.Range("N61").Value = -5353.5
$ws.Range("H113").Value = 6025.875
$ws.Range("J113").Value = 4949.5
$ws.Range("L113").Value = 4949.5
$ws.Range("N113").Value = -9289.5
$ws.Range("H117").Value = 57062.668
$ws.Range("J117").Value = 57062.668
$ws.Range("L117").Value = 57062.668
$ws.Range("N117").Value = -66240.66800000001
$ws.Range("H132").Value = 3262.75
$ws.Range("I132").Value = 3262.75
$ws.Range("K132").Value = 9788.25
$ws.Range("M132").Value = -7258.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3268.9285
$ws.Range("I113").Value = 4686.222
$ws.Range("J113").Value = 717.8
$ws.Range("K113").Value = 14058.666
$ws.Range("L113").Value = 2153.4
$ws.Range("M113").Value = -11888.666
$ws.Range("N113").Value = -6493.4

